# Refresh the cryptos price/volume table with the latest scraped values.
# Numeric-looking price strings (e.g. "303.23") must stay as TEXT (matching
# the source data, which always stores Price/Volume as text), so for those
# we force the Text number format before assigning the value and then clear
# the format again afterwards so the cell's style is left untouched (same
# as the original, unstyled cells).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.177.27"
$ws.Range("E2").Value = "  +1.81%  "
$ws.Range("D3").Value = "2.381.97"
$ws.Range("E3").Value = "  +4.14%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "303.23"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +0.75%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "96.96"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +1.68%  "
$ws.Range("E7").Value = "  +0.26%  "
$ws.Range("E8").Value = "  -0.08%  "
$ws.Range("E9").Value = "  +1.97%  "
$ws.Range("E10").Value = "  -0.73%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0790"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +1.32%  "
$ws.Range("E12").Value = "  +2.58%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "18.47"
$ws.Range("D13").ClearFormats()
$ws.Range("E14").Value = "  +0.90%  "
$ws.Range("D15").Value = "2.751.53"
$ws.Range("E15").Value = "  +4.04%  "
$ws.Range("D16").Value = "2.361.24"
$ws.Range("E16").Value = "  +3.44%  "
$ws.Range("E17").Value = "  +3.92%  "
$ws.Range("D18").Value = "43.156.50"
$ws.Range("E18").Value = "  +1.89%  "
$ws.Range("E19").Value = "  +0.60%  "
$ws.Range("E20").Value = "  +6.76%  "
$ws.Range("E21").Value = "  +0.44%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "68.53"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +1.34%  "
$ws.Range("E23").Value = "  +0.26%  "
$ws.Range("E24").Value = "  -2.20%  "
$ws.Range("B25").Value = "PancakeSwap"
$ws.Range("C25").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.44"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +1.36%  "
$ws.Range("B26").Value = "Dai"
$ws.Range("C26").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.00"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -0.08%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "24.83"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +2.35%  "
$ws.Range("E28").Value = "  +0.44%  "
$ws.Range("E29").Value = "  +1.20%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "31.58"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -0.56%  "
$ws.Range("E31").Value = "  +0.01%  "
$ws.Range("E32").Value = "  +2.62%  "
$ws.Range("E33").Value = "  +6.16%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "17.14"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -2.05%  "
$ws.Range("E35").Value = "  +7.17%  "
$ws.Range("E36").Value = "  +2.61%  "
$ws.Range("E37").Value = "  -0.98%  "
$ws.Range("E38").Value = "  -1.16%  "
$ws.Range("E39").Value = "  +4.72%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "22.36"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +11.56%  "
$ws.Range("E41").Value = "  +0.48%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "107.71"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -34.76%  "
$ws.Range("D43").Value = "1.955.95"
$ws.Range("E43").Value = "  -0.12%  "
$ws.Range("E44").Value = "  +0.92%  "
$ws.Range("E45").Value = "  +2.53%  "
$ws.Range("E46").Value = "  +0.69%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "9.24"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -10.69%  "
$ws.Range("B48").Value = "MultiversX"
$ws.Range("C48").Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "52.83"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -0.20%  "
$ws.Range("B49").Value = "Stacks"
$ws.Range("C49").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.51"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  +2.81%  "
$ws.Range("B50").Value = "BitcoinSV"
$ws.Range("C50").Value = "https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "72.10"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +1.90%  "
$ws.Range("B51").Value = "TrustWalletToken"
$ws.Range("C51").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.14"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +1.13%  "
